$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$cA8 = $ws.Range("A8")
$chA8 = $cA8.Characters(21, 2)
$chA8.Text = "28"

$cC9 = $ws.Range("C9")
# replace right-most date first so the left-hand offset (27) stays valid
$chC9b = $cC9.Characters(47, 8)
$chC9b.Text = "7/13/2025"
$chC9a = $cC9.Characters(27, 9)
$chC9a.Text = "7/7/2025"

# --- Column H width (Excel bestFit recalculation after the data changed) ---
$ws.Columns.Item(8).ColumnWidth = 7.433768

# --- Data table edits (rows 15-30) ---
# Row 15
$ws.Range("C15").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("N15").Value = -56.521739130434

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 101
$ws.Range("J16").Value = 104
$ws.Range("K16").Value = -2.884615384615
$ws.Range("L16").Value = 2.020202020202
$ws.Range("M16").Value = -45.108695652173
$ws.Range("N16").Value = -84.532924961715

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -22.727272727272
$ws.Range("I17").Value = 142
$ws.Range("J17").Value = 163
$ws.Range("K17").Value = -12.883435582822
$ws.Range("L17").Value = -15.976331360946
$ws.Range("M17").Value = 33.962264150943
$ws.Range("N17").Value = -49.645390070922

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 500
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 58.823529411764
$ws.Range("I18").Value = 139
$ws.Range("J18").Value = 181
$ws.Range("K18").Value = -23.204419889502
$ws.Range("L18").Value = 17.796610169491
$ws.Range("M18").Value = -33.492822966507
$ws.Range("N18").Value = -79.678362573099

# Row 19
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 14.285714285714
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -8.928571428571
$ws.Range("I19").Value = 368
$ws.Range("J19").Value = 387
$ws.Range("K19").Value = -4.909560723514
$ws.Range("L19").Value = -8.684863523573
$ws.Range("M19").Value = 57.93991416309
$ws.Range("N19").Value = 33.818181818181

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = 400
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 500
$ws.Range("I20").Value = 79
$ws.Range("J20").Value = 63
$ws.Range("K20").Value = 25.396825396825
$ws.Range("L20").Value = -7.058823529411
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -82.863340563991

# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = 23.636363636363
$ws.Range("I21").Value = 841
$ws.Range("J21").Value = 912
$ws.Range("K21").Value = -7.785087719298
$ws.Range("L21").Value = -4.971751412429
$ws.Range("M21").Value = 3.190184049079
$ws.Range("N21").Value = -64.796986186689

# Row 22
$ws.Range("C22").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = -100
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = 35.714285714285
$ws.Range("M22").Value = -17.391304347826

# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 85
$ws.Range("J23").Value = 94
$ws.Range("K23").Value = -9.574468085106
$ws.Range("L23").Value = -22.727272727272
$ws.Range("M23").Value = 18.055555555555

# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 82.35294117647
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 8.988764044943
$ws.Range("I24").Value = 605
$ws.Range("J24").Value = 592
$ws.Range("K24").Value = 2.195945945945
$ws.Range("L24").Value = 7.651245551601
$ws.Range("M24").Value = -3.968253968253

# Row 25
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 75
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -43.75
$ws.Range("I25").Value = 162
$ws.Range("J25").Value = 203
$ws.Range("K25").Value = -20.197044334975
$ws.Range("L25").Value = 118.918918918919

# Row 26
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 30
$ws.Range("F26").Value = 35
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 205
$ws.Range("J26").Value = 301
$ws.Range("K26").Value = -31.893687707641
$ws.Range("L26").Value = -20.233463035019
$ws.Range("M26").Value = -23.507462686567

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null
$ws.Range("H27").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$ws.Range("I27").Value = 11
$ws.Range("K27").Value = -21.428571428571
$ws.Range("L27").Value = -15.384615384615

# Row 28
$ws.Range("C28").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = -28.571428571428
$ws.Range("J28").Value = 32
$ws.Range("K28").Value = -3.125
$ws.Range("L28").Value = -13.888888888888

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4122) | Out-Null
$ws.Range("H29").Value = -100
$ws.Range("K14").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122) | Out-Null
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -71.428571428571
$ws.Range("N29").Value = -95.121951219512

# Row 30
$ws.Range("D30").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = -100
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Range("H30").Value = -100
$ws.Range("K14").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -66.666666666666
$ws.Range("N30").Value = -94.594594594594

$excel.CutCopyMode = 0